$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number must be
# kept as text (matching the source inlineStr cells), so force a text
# number format before assigning the value.

# Row 2
$ws.Range("D2").Value = "84.105.10"
$ws.Range("E2").Value = "  +5.52%  "

# Row 3
$ws.Range("D3").Value = "3.295.28"
$ws.Range("E3").Value = "  +2.70%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.37"
$ws.Range("E5").Value = "  +4.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "636.55"
$ws.Range("E6").Value = "  -0.37%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.325"
$ws.Range("E7").Value = "  +25.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  -1.83%  "

# Row 10
$ws.Range("D10").Value = "3.287.78"
$ws.Range("E10").Value = "  +2.37%  "

# Row 11
$ws.Range("E11").Value = "  -0.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000279"
$ws.Range("E12").Value = "  +4.95%  "

# Row 13
$ws.Range("E13").Value = "  -0.10%  "

# Row 14
$ws.Range("D14").Value = "3.884.04"
$ws.Range("E14").Value = "  +2.25%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.41"
$ws.Range("E15").Value = "  -0.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.64"
$ws.Range("E16").Value = "  +2.88%  "

# Row 17
$ws.Range("D17").Value = "84.261.31"
$ws.Range("E17").Value = "  +5.87%  "

# Row 18
$ws.Range("D18").Value = "3.280.65"
$ws.Range("E18").Value = "  +2.23%  "

# Row 19
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.51"
$ws.Range("E19").Value = "  -0.87%  "

# Row 20
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("E20").Value = "  +6.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "449.77"
$ws.Range("E21").Value = "  +0.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.16"
$ws.Range("E22").Value = "  -2.77%  "

# Row 23
$ws.Range("E23").Value = "  -0.86%  "

# Row 24
$ws.Range("E24").Value = "  +5.39%  "

# Row 25
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.31"
$ws.Range("E25").Value = "  +9.85%  "

# Row 26
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.34"
$ws.Range("E26").Value = "  +13.24%  "

# Row 27
$ws.Range("D27").Value = "3.455.72"
$ws.Range("E27").Value = "  +2.66%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "78.20"
$ws.Range("E28").Value = "  +0.50%  "

# Row 29
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("E30").Value = "  +0.90%  "

# Row 31
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.996"
$ws.Range("E32").Value = "  +0.12%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "569.93"
$ws.Range("E33").Value = "  +1.45%  "

# Row 34
$ws.Range("E34").Value = "  +26.14%  "

# Row 35
$ws.Range("E35").Value = "  -0.80%  "

# Row 36
$ws.Range("E36").Value = "  -1.68%  "

# Row 37
$ws.Range("E37").Value = "  -0.93%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.35"
$ws.Range("E38").Value = "  +0.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.23"
$ws.Range("E39").Value = "  +8.64%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.413"
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("E42").Value = "  +12.18%  "

# Row 43
$ws.Range("E43").Value = "  +3.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.05"
$ws.Range("E44").Value = "  +13.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "159.78"
$ws.Range("E45").Value = "  -2.36%  "

# Row 46
$ws.Range("E46").Value = "  +0.04%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "191.42"
$ws.Range("E47").Value = "  -1.01%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.13"
$ws.Range("E48").Value = "  +4.92%  "

# Row 49
$ws.Range("E49").Value = "  +0.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.782"
$ws.Range("E50").Value = "  -2.52%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.45"
$ws.Range("E51").Value = "  +1.86%  "
